# Update the "想去人数" (want-to-go count, column F) values that changed
# between the previous crawl and the latest crawl (commit "Update gh-pages
# to output generated at 456a3b4").
#
# Sheet "展览" (rId1 / sheet1.xml)
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value  = 5681   # was 5680
$wsExhibit.Range("F4").Value  = 82     # was 81
$wsExhibit.Range("F5").Value  = 8      # was 7
$wsExhibit.Range("F11").Value = 14     # was 13
$wsExhibit.Range("F12").Value = 89     # was 87
$wsExhibit.Range("F14").Value = 2400   # was 2398
$wsExhibit.Range("F15").Value = 422    # was 415

# Sheet "全部类型" (rId4 / sheet4.xml) mirrors the same events, but with
# the rows for the "演出" sheet interleaved, so the row numbers differ.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value  = 5681  # was 5680
$wsAll.Range("F5").Value  = 82    # was 81
$wsAll.Range("F6").Value  = 8     # was 7
$wsAll.Range("F13").Value = 14    # was 13
$wsAll.Range("F15").Value = 89    # was 87
$wsAll.Range("F17").Value = 2400  # was 2398
$wsAll.Range("F18").Value = 422   # was 415
